$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12's "Expected Result" (H12) previously held the "Failed..." text tied
# to the maximum-length test case. That whole max-length test case (row 14)
# is being removed, and row 12 (the placeholder test case) should instead
# report "Passed" - matching the style/format already used by row 11's H11.
$srcFmt = $ws.Range("H11")
$dstFmt = $ws.Range("H12")
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)
$dstFmt.Value = "Passed"

# Remove the two "TC_Pass_03" (minimum length) and "TC_Pass_04" (maximum
# length) password test-case rows entirely; remaining rows shift up.
$ws.Range("A13:H14").EntireRow.Delete()

# Update the active selection to reflect where the cursor ends up now that
# the sheet is shorter.
$ws.Range("A17").Select()
